# Update column G ("K") values on Sheet1 to reflect the recalculated
# strike counts (K) instead of the previous "Strike#" derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 3
    6  = 0
    7  = 0
    8  = 2
    9  = 1
    10 = 2
    11 = 1
    12 = 0
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
